$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = 'Elsevier eBooks'
$ws.Range("G5").Value = 'https://openalex.org/S4306463230'
$ws.Range("D6").Value = 'Widespread use of at-home COVID-19 tests hampers determination of community COVID-19 incidence. Using nationwide data available through the US National Wastewater Surveillance System, we examined the performance of two wastewater metrics in predicting high case and hospitalizations rates both before and after widespread use of at-home tests.We performed area under the receiver operating characteristic (ROC) curve analysis (AUC) for two wastewater metrics-viral concentration relative to the peak of January 2022 ("wastewater percentile") and 15-day percent change in SARS-CoV-2 ("percent change"). Dichotomized reported cases (≥ 200 or <200 cases per 100,000) and new hospitalizations (≥ 10 or <10 per 100,000) were our dependent variables, stratified by calendar quarter. Using logistic regression, we assessed the performance of combining wastewater metrics.Among 268 counties across 22 states, wastewater percentile detected high reported case and hospitalizations rates in the first quarter of 2022 (AUC 0.95 and 0.86 respectively) whereas the percent change did not (AUC 0.54 and 0.49 respectively). A wastewater percentile of 51% maximized sensitivity (0.93) and specificity (0.82) for detecting high case rates. A model inclusive of both metrics performed no better than using wastewater percentile alone. The predictive capability of wastewater percentile declined over time (AUC 0.84 and 0.72 for cases for second and third quarters of 2022).Nationwide, county wastewater levels above 51% relative to the historic peak predicted high COVID rates and hospitalization in the first quarter of 2022, but performed less well in subsequent quarters. Decline over time in predictive performance of this metric likely reflects underreporting of cases, reduced testing, and possibly lower virulence of infection due to vaccines and treatments.'
$ws.Range("AF6").Value = 'c("https://openalex.org/W4200329650", "https://openalex.org/W3127156785", "https://openalex.org/W3009669391", "https://openalex.org/W4205317059", "https://openalex.org/W3005417802", "https://openalex.org/W3084498529", "https://openalex.org/W3171943759", "https://openalex.org/W4226296940", "https://openalex.org/W3036314732", "https://openalex.org/W3134376730")'
$ws.Range("F8").Value = 'Oxford University Press eBooks'
$ws.Range("G8").Value = 'https://openalex.org/S4306463708'
$ws.Range("AF10").Value = 'c("https://openalex.org/W4376527619", "https://openalex.org/W2377074369", "https://openalex.org/W2285349705", "https://openalex.org/W3013339694", "https://openalex.org/W4289732247", "https://openalex.org/W3044411392", "https://openalex.org/W2952193828", "https://openalex.org/W3024928848", "https://openalex.org/W3120735633", "https://openalex.org/W4303987336")'
$ws.Range("AF11").Value = 'c("https://openalex.org/W392621609", "https://openalex.org/W3199894417", "https://openalex.org/W2537599345", "https://openalex.org/W3016102847", "https://openalex.org/W2049507953", "https://openalex.org/W2389938625", "https://openalex.org/W4386493544", "https://openalex.org/W1555246577", "https://openalex.org/W2979584651", "https://openalex.org/W2357077499")'
$ws.Range("AF12").Value = 'c("https://openalex.org/W2356237291", "https://openalex.org/W2960896737", "https://openalex.org/W2390208690", "https://openalex.org/W2371365383", "https://openalex.org/W2386180132", "https://openalex.org/W2383392149", "https://openalex.org/W2392428765", "https://openalex.org/W2355941487", "https://openalex.org/W3012488370", "https://openalex.org/W2112505184")'
